# required_subject_list.xlsx update:
#  - rename the single existing sheet
#  - add three new "Global SW" sheets, each populated by copying four rows
#    (COMP204 / COME331 / COMP319 / COMP312) out of the original table and
#    re-pointing their "교과구분" (and, for the COMP319 row, the course
#    code/name) at the new Global-SW values
#  - tweak sheet1's view

$wb = $excel.ActiveWorkbook

# ---- 1. rename the original sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "심화컴퓨터전공(ABEEK)"

# ---- 2. add the three new sheets, each after the previous one ---------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws2.Name = "글로벌소프트웨어전공(다중전공트랙)"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws3.Name = "글로벌소프트웨어전공(해외복수학위트랙)"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws4.Name = "글로벌소프트웨어전공(학석사연계트랙)"

# ---- 3. fill each new sheet with the Global SW course table -----------------------
foreach ($sheet in @($ws2, $ws3, $ws4)) {
    # header row, unchanged
    $ws1.Range("A1:F1").Copy()
    $sheet.Range("A1").PasteSpecial()

    # 프로그래밍기초 (COMP204) / 자료구조 (COME331)
    $ws1.Range("A7:F8").Copy()
    $sheet.Range("A2").PasteSpecial()

    # 알고리즘1 (COMP319) row -> becomes 알고리즘실습 (GLSO216)
    $ws1.Range("A12:F12").Copy()
    $sheet.Range("A4").PasteSpecial()

    # 운영체제 (COMP312)
    $ws1.Range("A13:F13").Copy()
    $sheet.Range("A5").PasteSpecial()

    $excel.CutCopyMode = $false

    # new Global-SW course: 알고리즘실습 / GLSO216
    $sheet.Cells.Item(4, 3).Value = "알고리즘실습"
    $sheet.Cells.Item(4, 1).Value = "GLSO216"

    # 교과구분 for every data row becomes 전공
    $sheet.Cells.Item(2, 4).Value = "전공"
    $sheet.Cells.Item(3, 4).Value = "전공"
    $sheet.Cells.Item(4, 4).Value = "전공"
    $sheet.Cells.Item(5, 4).Value = "전공"
}

# ---- 4. leave the 4th (newest) sheet active/selected -------------------------------
$ws4.Activate()

Write-Output "done"
